$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.509.05'
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = '1.868.25'
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '328.86'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = '0.4694'
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("D8").Value = '0.3975'
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("D9").Value = '47.01'
$ws.Range("E9").Value = '  -11.55%  '
$ws.Range("D10").Value = '0.08035'
$ws.Range("E10").Value = '  -4.14%  '
$ws.Range("D11").Value = '1.022'
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").Value = '21.74'
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '1.859.47'
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '86.95'
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").Value = '0.06566'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").Value = '17.22'
$ws.Range("E20").Value = '  -3.61%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '5.515'
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("D23").Value = '27.521.81'
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").Value = '2.308'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '2.089.82'
$ws.Range("E26").Value = '  -1.97%  '
$ws.Range("D27").Value = '20.35'
$ws.Range("E27").Value = '  +1.72%  '
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").Value = '2.085'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").Value = '5.564'
$ws.Range("E30").Value = '  -3.04%  '
$ws.Range("D31").Value = '122.56'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '0.09476'
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").Value = '0.9558'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").Value = '1.472'
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '3.604'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").Value = '5.305'
$ws.Range("E36").Value = '  -4.25%  '
$ws.Range("D37").Value = '0.06092'
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("D39").Value = '1.219'
$ws.Range("E39").Value = '  -4.81%  '
$ws.Range("D40").Value = '8.174'
$ws.Range("E40").Value = '  -6.87%  '
$ws.Range("D41").Value = '0.6002'
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '0.1898'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("E44").Value = '  -5.46%  '
$ws.Range("D45").Value = '1.267'
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("D46").Value = '0.5698'
$ws.Range("E46").Value = '  -2.70%  '
$ws.Range("D47").Value = '12.15'
$ws.Range("E47").Value = '  -5.19%  '
$ws.Range("D48").Value = '3.407'
$ws.Range("D49").Value = '1.939'
$ws.Range("E49").Value = '  -3.69%  '
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("D51").Value = '110.02'
$ws.Range("E51").Value = '  -0.79%  '
